# Weekly price update for Hortaliza, Macroferia Regional de Talca - Betarraga
# Inserts two new weekly records at the top of the data block (rows 469-470),
# pushing the existing historical rows (old 469-504) down to 471-506.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 469 (old first data row
# of this block). Inserting twice at the same index shifts everything below
# down by two rows total.
$ws.Rows.Item(469).Insert()
$ws.Rows.Item(469).Insert()

# --- New row 469 ---
$ws.Range("A469").Value = 5
$ws.Range("B469").Value = "Macroferia Regional de Talca"
$ws.Range("C469").Value = "Maule"
$ws.Range("D469").Value = 45013
$ws.Range("E469").Value = 7
$ws.Range("F469").Value = 100114014
$ws.Range("G469").Value = "Betarraga"
$ws.Range("H469").Value = "Sin especificar"
$ws.Range("I469").Value = "Primera"
$ws.Range("J469").Value = 4000
$ws.Range("K469").Value = 600
$ws.Range("L469").Value = 600
$ws.Range("M469").Value = 600
$ws.Range("N469").Value = "$/paquete 5 unidades"
$ws.Range("O469").Value = "Región del Maule"
$ws.Range("P469").Value = 120
$ws.Range("Q469").Value = 5
$ws.Range("R469").Value = "Hortaliza"

# --- New row 470 ---
$ws.Range("A470").Value = 5
$ws.Range("B470").Value = "Macroferia Regional de Talca"
$ws.Range("C470").Value = "Maule"
$ws.Range("D470").Value = 45013
$ws.Range("E470").Value = 7
$ws.Range("F470").Value = 100114014
$ws.Range("G470").Value = "Betarraga"
$ws.Range("H470").Value = "Sin especificar"
$ws.Range("I470").Value = "Segunda"
$ws.Range("J470").Value = 2000
$ws.Range("K470").Value = 500
$ws.Range("L470").Value = 500
$ws.Range("M470").Value = 500
$ws.Range("N470").Value = "$/paquete 5 unidades"
$ws.Range("O470").Value = "Región del Maule"
$ws.Range("P470").Value = 100
$ws.Range("Q470").Value = 5
$ws.Range("R470").Value = "Hortaliza"

# Ensure the style used for the date column (D) on the two new rows matches
# the rest of the data block (date/time number format).
$ws.Range("D469").NumberFormat = $ws.Range("D471").NumberFormat
$ws.Range("D470").NumberFormat = $ws.Range("D471").NumberFormat
